# Add team record (Wins / Losses / Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD, AE, AF ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the existing header formatting (bold, centered, bordered) from an
# adjacent header cell onto the three new header cells so they match the
# rest of row 1 (reuses the existing style instead of creating a new one).
$headerSrc = $ws.Range("AC1")
$headerSrc.Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (2-52): constant team record for every player row ---
for ($row = 2; $row -le 52; $row++) {
    $ws.Cells.Item($row, 30).Value = 74
    $ws.Cells.Item($row, 31).Value = 88
    $ws.Cells.Item($row, 32).Value = 0
}
